# Update gh-pages output data (想去人数 counts) in "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 272
$ws1.Range("F4").Value  = 6925
$ws1.Range("F7").Value  = 182
$ws1.Range("F9").Value  = 1141
$ws1.Range("F10").Value = 16495
$ws1.Range("F11").Value = 13
$ws1.Range("F13").Value = 56
$ws1.Range("F14").Value = 353
$ws1.Range("F17").Value = 11503
$ws1.Range("F18").Value = 21
$ws1.Range("F19").Value = 1151
$ws1.Range("F20").Value = 4545
$ws1.Range("F21").Value = 387
$ws1.Range("F22").Value = 395
$ws1.Range("F24").Value = 864
$ws1.Range("F25").Value = 327

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 272
$ws4.Range("F4").Value  = 6925
$ws4.Range("F7").Value  = 182
$ws4.Range("F10").Value = 1141
$ws4.Range("F11").Value = 16495
$ws4.Range("F12").Value = 13
$ws4.Range("F14").Value = 56
$ws4.Range("F15").Value = 353
$ws4.Range("F20").Value = 11504
$ws4.Range("F21").Value = 21
$ws4.Range("F22").Value = 1151
$ws4.Range("F23").Value = 4545
$ws4.Range("F24").Value = 387
$ws4.Range("F25").Value = 395
$ws4.Range("F27").Value = 864
$ws4.Range("F28").Value = 327
